$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "book" identifiers to "story" identifiers in the "Relates To ID" column (J)
$ws.Range("J3").Value = "SCH_001"
$ws.Range("J7").Value = "SCH_002"
$ws.Range("J10").Value = "SCH_003"

# Match the new selection left in the sheet (whole column J selected)
$ws.Range("J:J").Select()
